$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header columns for transformer parameters
$ws.Range("K1").Value = "d_model"
$ws.Range("L1").Value = "num_layers"

# Updated existing parameter values
$ws.Range("C2").Value = 23
$ws.Range("D2").Value = 90
$ws.Range("F2").Value = 0.0005
$ws.Range("H2").Value = 10

# New parameter values
$ws.Range("K2").Value = 16
$ws.Range("L2").Value = 2

# Update selection to match target state
$ws.Range("F5").Select()
